$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings are not
# auto-converted to numbers by Excel, matching the original inline string type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.124.71"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").Value = "1.828.84"
$ws.Range("E3").Value = "  -3.17%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "231.32"
$ws.Range("E5").Value = "  -2.80%  "

$ws.Range("D7").Value = "0.4658"
$ws.Range("E7").Value = "  -3.50%  "

$ws.Range("E8").Value = "  -6.58%  "

$ws.Range("D9").Value = "0.06261"
$ws.Range("E9").Value = "  -4.45%  "

$ws.Range("D10").Value = "1.823.31"
$ws.Range("E10").Value = "  -3.71%  "

$ws.Range("D11").Value = "0.07387"
$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("D12").Value = "15.98"
$ws.Range("E12").Value = "  -4.55%  "

$ws.Range("D13").Value = "4.899"
$ws.Range("E13").Value = "  -3.82%  "

$ws.Range("D14").Value = "83.15"
$ws.Range("E14").Value = "  -5.23%  "

$ws.Range("D15").Value = "0.6173"
$ws.Range("E15").Value = "  -7.42%  "

$ws.Range("D16").Value = "30.047.91"
$ws.Range("E16").Value = "  -1.79%  "

$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "226.03"
$ws.Range("E18").Value = "  -2.22%  "

$ws.Range("D19").Value = "0.000007262"
$ws.Range("E19").Value = "  -4.18%  "

$ws.Range("E20").Value = "  -6.31%  "

$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "2.072.35"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "4.840"
$ws.Range("E23").Value = "  -8.07%  "

$ws.Range("D24").Value = "5.863"
$ws.Range("E24").Value = "  -5.28%  "

$ws.Range("D25").Value = "9.126"
$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("D26").Value = "164.55"
$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("D27").Value = "17.59"
$ws.Range("E27").Value = "  -5.84%  "

$ws.Range("D28").Value = "1.846"
$ws.Range("E28").Value = "  -5.56%  "

$ws.Range("D29").Value = "0.1014"
$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("E30").Value = "  -2.30%  "

$ws.Range("D31").Value = "4.046"
$ws.Range("E31").Value = "  -6.50%  "

$ws.Range("D32").Value = "3.761"
$ws.Range("E32").Value = "  -6.42%  "

$ws.Range("D33").Value = "0.04777"
$ws.Range("E33").Value = "  -5.37%  "

$ws.Range("D34").Value = "1.122"
$ws.Range("E34").Value = "  -6.86%  "

$ws.Range("D35").Value = "0.7035"
$ws.Range("E35").Value = "  -6.16%  "

$ws.Range("D36").Value = "2.687"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").Value = "0.01813"
$ws.Range("E37").Value = "  -3.58%  "

$ws.Range("D38").Value = "2.600"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").Value = "0.8926"
$ws.Range("E39").Value = "  -3.01%  "

$ws.Range("D40").Value = "1.923"
$ws.Range("E40").Value = "  -6.78%  "

$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").Value = "103.28"
$ws.Range("E42").Value = "  -3.52%  "

$ws.Range("D43").Value = "5.468"
$ws.Range("E43").Value = "  -3.22%  "

$ws.Range("D44").Value = "0.4001"
$ws.Range("E44").Value = "  -6.62%  "

$ws.Range("D45").Value = "6.969"
$ws.Range("E45").Value = "  -6.13%  "

$ws.Range("D46").Value = "0.1189"
$ws.Range("E46").Value = "  -6.84%  "

$ws.Range("D47").Value = "59.61"
$ws.Range("E47").Value = "  -7.01%  "

$ws.Range("D48").Value = "8.428"
$ws.Range("E48").Value = "  -6.30%  "

$ws.Range("D49").Value = "32.62"
$ws.Range("E49").Value = "  -4.11%  "

$ws.Range("D50").Value = "0.05514"
$ws.Range("E50").Value = "  -2.71%  "

$ws.Range("D51").Value = "1.367"
$ws.Range("E51").Value = "  -7.88%  "

# Restore default (Normal) style for column D so no residual number-format
# style is left applied to the cells (matches original unstyled cells).
$ws.Range("D2:D51").Style = "Normal"
